# Natmi following Dr Hou advice
#
# The LR-pairs table (Sema3a -> Nrp1) is recomputed: "Sending cluster" /
# "Target cluster" now range over 3 cell populations (ECs, FAPs, sCs)
# instead of 2 (FAPs, sCs), giving a 3x3 = 9 row table (rows 2-10) instead
# of the previous 2x3 = 6 row table (rows 2-7), and all the numeric
# columns (E:T) are refreshed with the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Sema3a/Nrp1 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3a"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5674196666666667
$ws.Range("H2").Value = 1.702259
$ws.Range("I2").Value = 0.07864125446886469
$ws.Range("J2").Value = 0.07864125446886468
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 102.8289443333334
$ws.Range("N2").Value = 308.486833
$ws.Range("O2").Value = 0.5559120396302444
$ws.Range("P2").Value = 0.5559120396302443
$ws.Range("Q2").Value = 58.34716531730523
$ws.Range("R2").Value = 525.124487855747
$ws.Range("S2").Value = 0.04371762017086764
$ws.Range("T2").Value = 0.04371762017086762

# Row 3: ECs -> Sema3a/Nrp1 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3a"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5674196666666667
$ws.Range("H3").Value = 1.702259
$ws.Range("I3").Value = 0.07864125446886469
$ws.Range("J3").Value = 0.07864125446886468
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3441717873742006
$ws.Range("P3").Value = 0.3441717873742006
$ws.Range("Q3").Value = 36.12342734802389
$ws.Range("R3").Value = 325.110846132215
$ws.Range("S3").Value = 0.0270661011118985
$ws.Range("T3").Value = 0.0270661011118985

# Row 4: ECs -> Sema3a/Nrp1 -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3a"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5674196666666667
$ws.Range("H4").Value = 1.702259
$ws.Range("I4").Value = 0.07864125446886469
$ws.Range("J4").Value = 0.07864125446886468
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.481835
$ws.Range("N4").Value = 55.445505
$ws.Range("O4").Value = 0.09991617299555507
$ws.Range("P4").Value = 0.09991617299555505
$ws.Range("Q4").Value = 10.48695665508833
$ws.Range("R4").Value = 94.382609895795
$ws.Range("S4").Value = 0.007857533186098553
$ws.Range("T4").Value = 0.00785753318609855

# Row 5: FAPs -> Sema3a/Nrp1 -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema3a"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7227763333333334
$ws.Range("H5").Value = 2.168329
$ws.Range("I5").Value = 0.1001728366019618
$ws.Range("J5").Value = 0.1001728366019618
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 102.8289443333334
$ws.Range("N5").Value = 308.486833
$ws.Range("O5").Value = 0.5559120396302444
$ws.Range("P5").Value = 0.5559120396302443
$ws.Range("Q5").Value = 74.32232734578413
$ws.Range("R5").Value = 668.9009461120571
$ws.Range("S5").Value = 0.05568728591094378
$ws.Range("T5").Value = 0.05568728591094377

# Row 6: FAPs -> Sema3a/Nrp1 -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3a"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7227763333333334
$ws.Range("H6").Value = 2.168329
$ws.Range("I6").Value = 0.1001728366019618
$ws.Range("J6").Value = 0.1001728366019618
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3441717873742006
$ws.Range("P6").Value = 0.3441717873742006
$ws.Range("Q6").Value = 46.01384107712945
$ws.Range("R6").Value = 414.124569694165
$ws.Range("S6").Value = 0.03447666421964094
$ws.Range("T6").Value = 0.03447666421964094

# Row 7: FAPs -> Sema3a/Nrp1 -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3a"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7227763333333334
$ws.Range("H7").Value = 2.168329
$ws.Range("I7").Value = 0.1001728366019618
$ws.Range("J7").Value = 0.1001728366019618
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.481835
$ws.Range("N7").Value = 55.445505
$ws.Range("O7").Value = 0.09991617299555507
$ws.Range("P7").Value = 0.09991617299555505
$ws.Range("Q7").Value = 13.35823293457167
$ws.Range("R7").Value = 120.224096411145
$ws.Range("S7").Value = 0.01000888647137709
$ws.Range("T7").Value = 0.01000888647137708

# Row 8: sCs -> Sema3a/Nrp1 -> ECs (new row)
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema3a"
$ws.Range("C8").Value = "Nrp1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.925096666666668
$ws.Range("H8").Value = 17.77529
$ws.Range("I8").Value = 0.8211859089291735
$ws.Range("J8").Value = 0.8211859089291734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 102.8289443333334
$ws.Range("N8").Value = 308.486833
$ws.Range("O8").Value = 0.5559120396302444
$ws.Range("P8").Value = 0.5559120396302443
$ws.Range("Q8").Value = 609.2714353062858
$ws.Range("R8").Value = 5483.442917756572
$ws.Range("S8").Value = 0.456507133548433
$ws.Range("T8").Value = 0.4565071335484328

# Row 9: sCs -> Sema3a/Nrp1 -> FAPs (new row)
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema3a"
$ws.Range("C9").Value = "Nrp1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.925096666666668
$ws.Range("H9").Value = 17.77529
$ws.Range("I9").Value = 0.8211859089291735
$ws.Range("J9").Value = 0.8211859089291734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 63.66262833333334
$ws.Range("N9").Value = 190.987885
$ws.Range("O9").Value = 0.3441717873742006
$ws.Range("P9").Value = 0.3441717873742006
$ws.Range("Q9").Value = 377.2072269290723
$ws.Range("R9").Value = 3394.86504236165
$ws.Range("S9").Value = 0.2826290220426612
$ws.Range("T9").Value = 0.2826290220426612

# Row 10: sCs -> Sema3a/Nrp1 -> sCs (new row)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema3a"
$ws.Range("C10").Value = "Nrp1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.925096666666668
$ws.Range("H10").Value = 17.77529
$ws.Range("I10").Value = 0.8211859089291735
$ws.Range("J10").Value = 0.8211859089291734
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 18.481835
$ws.Range("N10").Value = 55.445505
$ws.Range("O10").Value = 0.09991617299555507
$ws.Range("P10").Value = 0.09991617299555505
$ws.Range("Q10").Value = 109.5066589523833
$ws.Range("R10").Value = 985.5599305714501
$ws.Range("S10").Value = 0.08204975333807943
$ws.Range("T10").Value = 0.08204975333807941
